# Update gh-pages to output generated at 7921097
# Increments the "想去人数" (number of people wanting to go) column F
# across all four worksheets to reflect freshly scraped counts.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (Exhibitions)
$wsExhibit.Range("F2").Value  = 11420
$wsExhibit.Range("F3").Value  = 1941
$wsExhibit.Range("F4").Value  = 529
$wsExhibit.Range("F11").Value = 1349
$wsExhibit.Range("F14").Value = 14
$wsExhibit.Range("F16").Value = 540
$wsExhibit.Range("F17").Value = 678
$wsExhibit.Range("F20").Value = 938
$wsExhibit.Range("F21").Value = 9
$wsExhibit.Range("F28").Value = 681

# Sheet "演出" (Shows)
$wsShow.Range("F7").Value  = 8
$wsShow.Range("F8").Value  = 100
$wsShow.Range("F10").Value = 391

# Sheet "本地生活" (Local life)
$wsLocal.Range("F2").Value = 51

# Sheet "全部类型" (All types)
$wsAll.Range("F2").Value  = 11420
$wsAll.Range("F3").Value  = 1941
$wsAll.Range("F5").Value  = 529
$wsAll.Range("F13").Value = 51
$wsAll.Range("F14").Value = 1349
$wsAll.Range("F19").Value = 14
$wsAll.Range("F21").Value = 540
$wsAll.Range("F22").Value = 678
$wsAll.Range("F25").Value = 938
$wsAll.Range("F26").Value = 9
$wsAll.Range("F32").Value = 8
$wsAll.Range("F33").Value = 100
$wsAll.Range("F34").Value = 100
$wsAll.Range("F37").Value = 681
$wsAll.Range("F41").Value = 391
